$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet title/name and the "through" date in the header row
$ws.Name = "Through 2022-10-05"
$ws.Range("A11").Value = "October (through 10-05)"

# Update October row (row 11)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 12
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 29
$ws.Range("H11").Value = 35
$ws.Range("I11").Value = 17

# Update Total row (row 12)
$ws.Range("B12").Value = 231
$ws.Range("C12").Value = 436
$ws.Range("D12").Value = 638
$ws.Range("E12").Value = 560
$ws.Range("F12").Value = 426
$ws.Range("G12").Value = 930
$ws.Range("H12").Value = 1282
$ws.Range("I12").Value = 1299
